$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.043.03"
$ws.Range("E2").Value = "  -4.32%  "
$ws.Range("D3").Value = "2.234.57"
$ws.Range("E3").Value = "  -5.01%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.90"
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  -6.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.53"
$ws.Range("E7").Value = "  -5.59%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -6.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0976"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.21"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "34.98"
$ws.Range("E12").Value = "  +7.02%  "
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.68"
$ws.Range("E14").Value = "  -7.14%  "
$ws.Range("D15").Value = "2.574.23"
$ws.Range("E15").Value = "  -4.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.89"
$ws.Range("E16").Value = "  -8.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("E17").Value = "  -6.06%  "
$ws.Range("D18").Value = "2.232.46"
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("D19").Value = "42.006.80"
$ws.Range("E19").Value = "  -4.26%  "
$ws.Range("E20").Value = "  -5.74%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.23"
$ws.Range("E21").Value = "  -7.31%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.02"
$ws.Range("E22").Value = "  -6.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.09"
$ws.Range("E23").Value = "  -7.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.34"
$ws.Range("E27").Value = "  -6.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.96"
$ws.Range("E28").Value = "  -4.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -4.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.67"
$ws.Range("E30").Value = "  -4.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.44"
$ws.Range("E31").Value = "  -8.97%  "
$ws.Range("E32").Value = "  -6.77%  "
$ws.Range("E33").Value = "  -7.43%  "
$ws.Range("E34").Value = "  -4.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("E36").Value = "  -8.67%  "
$ws.Range("E37").Value = "  -4.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.45"
$ws.Range("E38").Value = "  +13.99%  "
$ws.Range("E39").Value = "  -5.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.95"
$ws.Range("E40").Value = "  -7.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0263"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.59"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.84"
$ws.Range("E43").Value = "  -4.61%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.76"
$ws.Range("E44").Value = "  -6.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.101"
$ws.Range("E45").Value = "  -6.80%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.187"
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("B48").Value = "BitTorrent-New"
$ws.Range("C48").Value = "https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt"
$ws.Range("D48").Value = "0.0₃0154"
$ws.Range("E48").Value = "  +14.81%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -5.41%  "
$ws.Range("E50").Value = "  -6.26%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.26"
$ws.Range("E51").Value = "  +5.56%  "
